# Applies the LOM3037.xlsx content update:
#  - removes the two "docentes" rows (5840963 - Daniela Camargo Vernilli /
#    1922320 - Sebastiao Ribeiro) that used to sit right under "Objetivos:",
#    shifting every row below them up by two
#  - clears out the long free-text paragraphs (Objetivos, Programa resumido,
#    Programa, Bibliografia) that used to live in columns B/C, replacing
#    them with the (stray/duplicated) short values that remain in the
#    published sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the two rows that held the "Daniela" / "Sebastiao" docente
#    values right after "Objetivos:" - this shifts rows 15-25 up to 13-23
#    and keeps every row's height/labels in column A intact automatically.
$ws.Rows("13:14").Delete() | Out-Null

# 2) Column B/C content fix-ups for the rows whose long paragraph text was
#    removed. Most of these are plain literal strings, safe to type
#    directly. "01/01/2020" is date-like, so typing it naively would make
#    Excel auto-convert it to a date serial / new number format; instead
#    we copy the *value* of the existing "01/01/2020" text cell (B8/C8)
#    so the destination keeps its original (text) type and style.

# Row 10 (Objetivos:) - B/C used to hold the long objectives paragraph.
$ws.Range("B10").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("C10").Value = "5840963 - Daniela Camargo Vernilli"

# Row 13 (Programa resumido:) - B/C used to hold the short-syllabus text.
# xlPasteValues = -4163
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4163) | Out-Null
$ws.Range("C8").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Row 15 (Programa:) - B/C used to hold the long syllabus paragraph.
$ws.Range("B15").Value = "5840963 - Daniela Camargo Vernilli"
$ws.Range("C15").Value = "5840963 - Daniela Camargo Vernilli"

# Row 18 (Método:) - B/C used to hold the "Duas provas..." text.
$ws.Range("B18").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C18").Value = "1922320 - Sebastiao Ribeiro"

# Row 19 (Critério:) - B/C used to hold the "MS= P1+P2/2..." text.
$ws.Range("B19").Value = "Duas provas bimestrais escritas (P1 e P2), cada uma valendo nota de 0,0  a 10,0."
$ws.Range("C19").Value = "Duas provas bimestrais escritas (P1 e P2), cada uma valendo nota de 0,0  a 10,0."

# Row 20 (Norma de recuperação:) - B/C used to hold the "Estudo dirigido..." text.
$ws.Range("B20").Value = "MS= P1+P2/2, onde: MS= média do semestre.MS> ou = 5,0 = Aluno AprovadoMS< 3,0 = Aluno Reprovado3,0 < ou = MS < 5,0 = Aluno de Recuperação."
$ws.Range("C20").Value = "MS= P1+P2/2, onde: MS= média do semestre.MS> ou = 5,0 = Aluno AprovadoMS< 3,0 = Aluno Reprovado3,0 < ou = MS < 5,0 = Aluno de Recuperação."

# Row 21 (Bibliografia:) - B/C used to hold the long bibliography paragraph.
$ws.Range("B21").Value = "Estudo dirigido de todo o conteúdo da disciplina e uma prova (PR) valendo nota de 0,0 a 10,0, contendo todo o conteúdo da disciplina.O aluno será aprovado se apresentar (média final) MF > ou = 5,0.Onde: MF= MS+PR/2, onde:  MS= média do semestre e PR= prova de recuperação."
$ws.Range("C21").Value = "Estudo dirigido de todo o conteúdo da disciplina e uma prova (PR) valendo nota de 0,0 a 10,0, contendo todo o conteúdo da disciplina.O aluno será aprovado se apresentar (média final) MF > ou = 5,0.Onde: MF= MS+PR/2, onde:  MS= média do semestre e PR= prova de recuperação."

Write-Host "Edit applied"
